# Auto-generated script to update cryptos list (Price/Volume columns)
# per commit: 'Updated cryptos list on Sat Feb 25 14:56:58 UTC 2023 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.045.02'
$ws.Range("E2").Value = '  -3.38%  '

$ws.Range("D3").Value = '1.599.91'
$ws.Range("E3").Value = '  -2.38%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.96'
$ws.Range("E6").Value = '  -2.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3785'
$ws.Range("E7").Value = '  -1.93%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3648'
$ws.Range("E8").Value = '  -4.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '50.01'
$ws.Range("E9").Value = '  -1.35%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.261'
$ws.Range("E10").Value = '  -4.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").Value = '  +0.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08143'
$ws.Range("E12").Value = '  -2.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.59'
$ws.Range("E13").Value = '  -4.50%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.598'
$ws.Range("E14").Value = '  -4.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.361'
$ws.Range("E15").Value = '  -5.07%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001247'
$ws.Range("E16").Value = '  -4.15%  '

$ws.Range("D17").Value = '1.607.25'
$ws.Range("E17").Value = '  -3.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.08'
$ws.Range("E18").Value = '  -1.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06837'
$ws.Range("E19").Value = '  -1.33%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.26'
$ws.Range("E20").Value = '  -5.61%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.546'
$ws.Range("E21").Value = '  -4.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.5587'
$ws.Range("E22").Value = '  -5.49%  '

$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.07'
$ws.Range("E24").Value = '  -2.96%  '

$ws.Range("D25").Value = '23.024.28'
$ws.Range("E25").Value = '  -3.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.362'
$ws.Range("E26").Value = '  -2.93%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.829'
$ws.Range("E27").Value = '  -1.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.06'
$ws.Range("E28").Value = '  -3.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '150.35'
$ws.Range("E29").Value = '  -1.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.233'
$ws.Range("E30").Value = '  -4.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.21'
$ws.Range("E31").Value = '  -1.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.331'
$ws.Range("E32").Value = '  -5.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.826'
$ws.Range("E33").Value = '  -12.09%  '

$ws.Range("D34").Value = '1.783.45'
$ws.Range("E34").Value = '  -1.99%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9670'
$ws.Range("E35").Value = '  -1.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.07580'
$ws.Range("E36").Value = '  -4.34%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.34'
$ws.Range("E37").Value = '  -0.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.268'
$ws.Range("E38").Value = '  -4.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02706'
$ws.Range("E39").Value = '  -5.94%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2531'
$ws.Range("E40").Value = '  -4.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.08875'
$ws.Range("E41").Value = '  -2.00%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.365'
$ws.Range("E42").Value = '  -3.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7034'
$ws.Range("E43").Value = '  -5.75%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.40'
$ws.Range("E44").Value = '  -6.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.26'
$ws.Range("E45").Value = '  -7.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6629'
$ws.Range("E46").Value = '  -3.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9998'
$ws.Range("E47").Value = '  +0.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.291'
$ws.Range("E48").Value = '  -4.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.994'
$ws.Range("E49").Value = '  -1.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.65'

$ws.Range("E51").Value = '  -3.82%  '
